# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted as row 14 (Fecha 2023-01-25 / serial
# 44951), pushing the previously existing rows 14-17 down to rows 15-18
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 14, shifting rows 14:17 down to 15:18.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C14").Value = "Ñuble"
$ws.Range("D14").Value = 44951
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100101
$ws.Range("H14").Value = "Berries"
$ws.Range("I14").Value = 100101001
$ws.Range("J14").Value = "Arándano (blue)"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 2800
$ws.Range("O14").Value = 3000
$ws.Range("P14").Value = 2900
$ws.Range("Q14").Value = "$/bandeja 2 kilos"
$ws.Range("R14").Value = "Provincia de Diguillín"
$ws.Range("S14").Value = 1450
$ws.Range("T14").Value = 2
